# This script reproduces the commit "parses immediates, range not taken
# in account yet" on the Instructionset workbook.
#
# The substantive change is on the "Operands" worksheet: the old
# "Constant" header was renamed to "Key", two new columns "Min" and "Max"
# were appended, and six new rows describing the immediate operand types
# (IMM6, IMM6s, IMM8, IMM8s, IMM12, IMM12s) were added below the existing
# data, each carrying a bit width ("k") and a Min/Max range.
#
# It also updates the selection left behind on a few sheets after the
# editing session, with "Operands" becoming the active tab.

$wb = $excel.ActiveWorkbook

$wsMicro = $wb.Worksheets.Item("MicroInstructions")
$wsEnc   = $wb.Worksheets.Item("Encoding")
$wsOps   = $wb.Worksheets.Item("Operands")
$wsSum   = $wb.Worksheets.Item("Summary")

# --- Operands sheet: the real content edit -----------------------------
$wsOps.Activate()

# The author first typed out the IMM8 row (row 20)...
$wsOps.Cells.Item(20, 1).Value = "IMM8"
$wsOps.Cells.Item(20, 2).Value = "k"
$wsOps.Cells.Item(20, 3).Value = 8
$wsOps.Cells.Item(20, 6).Value = 0
$wsOps.Cells.Item(20, 7).Value = 255

# ...then renamed the old "Constant" header to "Key" and added the new
# "Min"/"Max" headers...
$wsOps.Cells.Item(1, 4).Value = "Key"
$wsOps.Cells.Item(1, 6).Value = "Min"
$wsOps.Cells.Item(1, 7).Value = "Max"

# ...and finally filled in the rest of the immediate rows.
$wsOps.Cells.Item(18, 1).Value = "IMM6"
$wsOps.Cells.Item(18, 2).Value = "k"
$wsOps.Cells.Item(18, 3).Value = 6
$wsOps.Cells.Item(18, 6).Value = 0
$wsOps.Cells.Item(18, 7).Value = 64

$wsOps.Cells.Item(19, 1).Value = "IMM6s"
$wsOps.Cells.Item(19, 2).Value = "k"
$wsOps.Cells.Item(19, 3).Value = 6
$wsOps.Cells.Item(19, 6).Value = -31
$wsOps.Cells.Item(19, 7).Value = 32

$wsOps.Cells.Item(21, 1).Value = "IMM8s"
$wsOps.Cells.Item(21, 2).Value = "k"
$wsOps.Cells.Item(21, 3).Value = 8
$wsOps.Cells.Item(21, 6).Value = -127
$wsOps.Cells.Item(21, 7).Value = 128

$wsOps.Cells.Item(22, 1).Value = "IMM12"
$wsOps.Cells.Item(22, 2).Value = "k"
$wsOps.Cells.Item(22, 3).Value = 12
$wsOps.Cells.Item(22, 6).Value = 0
$wsOps.Cells.Item(22, 7).Value = 4096

$wsOps.Cells.Item(23, 1).Value = "IMM12s"
$wsOps.Cells.Item(23, 2).Value = "k"
$wsOps.Cells.Item(23, 3).Value = 12
$wsOps.Cells.Item(23, 6).Value = -2047
$wsOps.Cells.Item(23, 7).Value = 2048

# --- Leftover cursor/selection state from the editing session ----------
$wsMicro.Activate()
$wsMicro.Range("A28").Select() | Out-Null

$wsEnc.Activate()
$wsEnc.Range("A7").Select() | Out-Null

$wsOps.Activate()
$wsOps.Range("F24").Select() | Out-Null
